# Meta Descriptions and UI updated for Adsense Ad
# Marks column F (TRUE/FALSE checkbox) as TRUE for the rows whose pages
# have had their Meta Description / Adsense UI work completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,73,77,84)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = $true
}

# Reflect where the user's selection/scroll position ended up after doing
# this work (scrolled down the list, cursor resting on F84).
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 69
    $win.ScrollColumn = 1
} catch {
    # Scroll-position isn't modeled in every host; ignore if unsupported.
}
$ws.Range("F84").Select()
